$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.430.66'
$ws.Range("E2").Value = '  +1.60%  '
$ws.Range("D3").Value = '2.280.95'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '307.23'
$ws.Range("E5").Value = '  +1.27%  '
$ws.Range("D6").Value = '97.67'
$ws.Range("E6").Value = '  +5.89%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +2.34%  '
$ws.Range("D10").Value = '35.62'
$ws.Range("E10").Value = '  +9.81%  '
$ws.Range("D11").Value = '0.0798'
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("D13").Value = '''6.70'
$ws.Range("E13").Value = '  +0.47%  '
$ws.Range("D14").Value = '2.631.32'
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").Value = '''14.50'
$ws.Range("E15").Value = '  +1.80%  '
$ws.Range("D16").Value = '2.278.22'
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '0.798'
$ws.Range("E17").Value = '  +3.66%  '
$ws.Range("D18").Value = '42.316.94'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").Value = '12.61'
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("D21").Value = '5.98'
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").Value = '67.71'
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("D23").Value = '240.77'
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").Value = '''2.60'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").Value = '1.95'
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '23.87'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").Value = '37.75'
$ws.Range("E28").Value = '  +6.97%  '
$ws.Range("D29").Value = '9.51'
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("E30").Value = '  +1.35%  '
$ws.Range("D31").Value = '159.82'
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").Value = '5.27'
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = '3.14'
$ws.Range("E34").Value = '  +4.20%  '
$ws.Range("D35").Value = '0.0742'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").Value = '17.07'
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("E38").Value = '  +1.21%  '
$ws.Range("D39").Value = '1.85'
$ws.Range("E39").Value = '  +3.23%  '
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("D41").Value = '4.13'
$ws.Range("E41").Value = '  +5.79%  '
$ws.Range("E42").Value = '  +14.19%  '
$ws.Range("D43").Value = '2.001.02'
$ws.Range("E43").Value = '  -0.38%  '
$ws.Range("E44").Value = '  +2.39%  '
$ws.Range("D45").Value = '18.96'
$ws.Range("E45").Value = '  -1.54%  '
$ws.Range("D46").Value = '2.99'
$ws.Range("E46").Value = '  +3.67%  '
$ws.Range("D47").Value = '10.02'
$ws.Range("E47").Value = '  -3.19%  '
$ws.Range("D48").Value = '53.01'
$ws.Range("E48").Value = '  +1.31%  '
$ws.Range("D49").Value = '1.53'
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("D50").Value = '''72.20'
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("D51").Value = '91.97'
$ws.Range("E51").Value = '  +1.04%  '
